$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.3915976666666667
$ws.Range("H2").Value = 1.174793
$ws.Range("I2").Value = 0.02606065131430495
$ws.Range("J2").Value = 0.02606065131430495
$ws.Range("M2").Value = 135.0916853333333
$ws.Range("N2").Value = 405.2750559999999
$ws.Range("O2").Value = 0.7123704212620513
$ws.Range("P2").Value = 0.7123704212620514
$ws.Range("Q2").Value = 52.90158876260089
$ws.Range("R2").Value = 476.114298863408
$ws.Range("S2").Value = 0.01856483715513485
$ws.Range("T2").Value = 0.01856483715513485

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.3915976666666667
$ws.Range("H3").Value = 1.174793
$ws.Range("I3").Value = 0.02606065131430495
$ws.Range("J3").Value = 0.02606065131430495
$ws.Range("O3").Value = 0.2125756143240238
$ws.Range("P3").Value = 0.2125756143240238
$ws.Range("Q3").Value = 15.78615197133511
$ws.Range("R3").Value = 142.075367742016
$ws.Range("S3").Value = 0.005539858962822553
$ws.Range("T3").Value = 0.005539858962822554

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.3915976666666667
$ws.Range("H4").Value = 1.174793
$ws.Range("I4").Value = 0.02606065131430495
$ws.Range("J4").Value = 0.02606065131430495
$ws.Range("M4").Value = 14.23299766666667
$ws.Range("N4").Value = 42.698993
$ws.Range("O4").Value = 0.07505396441392481
$ws.Range("P4").Value = 0.07505396441392483
$ws.Range("Q4").Value = 5.573608675938779
$ws.Range("R4").Value = 50.16247808344901
$ws.Range("S4").Value = 0.001955955196347547
$ws.Range("T4").Value = 0.001955955196347547

# Row 5
$ws.Range("I5").Value = 0.4187506438669658
$ws.Range("J5").Value = 0.4187506438669658
$ws.Range("M5").Value = 135.0916853333333
$ws.Range("N5").Value = 405.2750559999999
$ws.Range("O5").Value = 0.7123704212620513
$ws.Range("P5").Value = 0.7123704212620514
$ws.Range("Q5").Value = 850.0391678148426
$ws.Range("R5").Value = 7650.352510333583
$ws.Range("S5").Value = 0.2983055725752656
$ws.Range("T5").Value = 0.2983055725752657

# Row 6
$ws.Range("I6").Value = 0.4187506438669658
$ws.Range("J6").Value = 0.4187506438669658
$ws.Range("O6").Value = 0.2125756143240238
$ws.Range("P6").Value = 0.2125756143240238
$ws.Range("S6").Value = 0.08901617536860075
$ws.Range("T6").Value = 0.08901617536860076

# Row 7
$ws.Range("I7").Value = 0.4187506438669658
$ws.Range("J7").Value = 0.4187506438669658
$ws.Range("M7").Value = 14.23299766666667
$ws.Range("N7").Value = 42.698993
$ws.Range("O7").Value = 0.07505396441392481
$ws.Range("P7").Value = 0.07505396441392483
$ws.Range("Q7").Value = 89.55847624693634
$ws.Range("R7").Value = 806.0262862224271
$ws.Range("S7").Value = 0.03142889592309935
$ws.Range("T7").Value = 0.03142889592309936

# Row 8
$ws.Range("G8").Value = 8.342485333333334
$ws.Range("H8").Value = 25.027456
$ws.Range("I8").Value = 0.5551887048187292
$ws.Range("J8").Value = 0.5551887048187292
$ws.Range("M8").Value = 135.0916853333333
$ws.Range("N8").Value = 405.2750559999999
$ws.Range("O8").Value = 0.7123704212620513
$ws.Range("P8").Value = 0.7123704212620514
$ws.Range("Q8").Value = 1127.000403548615
$ws.Range("R8").Value = 10143.00363193754
$ws.Range("S8").Value = 0.3955000115316508
$ws.Range("T8").Value = 0.3955000115316508

# Row 9
$ws.Range("G9").Value = 8.342485333333334
$ws.Range("H9").Value = 25.027456
$ws.Range("I9").Value = 0.5551887048187292
$ws.Range("J9").Value = 0.5551887048187292
$ws.Range("O9").Value = 0.2125756143240238
$ws.Range("P9").Value = 0.2125756143240238
$ws.Range("Q9").Value = 336.3036925414969
$ws.Range("R9").Value = 3026.733232873472
$ws.Range("S9").Value = 0.1180195799926005
$ws.Range("T9").Value = 0.1180195799926005

# Row 10
$ws.Range("G10").Value = 8.342485333333334
$ws.Range("H10").Value = 25.027456
$ws.Range("I10").Value = 0.5551887048187292
$ws.Range("J10").Value = 0.5551887048187292
$ws.Range("M10").Value = 14.23299766666667
$ws.Range("N10").Value = 42.698993
$ws.Range("O10").Value = 0.07505396441392481
$ws.Range("P10").Value = 0.07505396441392483
$ws.Range("Q10").Value = 118.7385742835342
$ws.Range("R10").Value = 1068.647168551808
$ws.Range("S10").Value = 0.04166911329447791
$ws.Range("T10").Value = 0.04166911329447791
